# Projektkennzahlen.xlsx - update "6. Projektstatusbericht" data (row 8) and
# add data for the new "7. Projektstatusbericht" (row 9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (6. Projektstatusbericht) : SOLL/IST figures changed ---
# A8 used to be the formula =1248.94+(2579.69/5); it is now a plain input value.
$ws.Range("A8").Value = 3025.51
$ws.Range("B8").Value = 2004.75
$ws.Range("E8").Value = 2711.34
$ws.Range("F8").Value = 2004.75

# C8/G8/I8/J8/K8 keep their existing formulas (=A8+B8, =E8+F8, =A8-E8, =B8-F8, =C8-G8)
# and simply recalculate with the new inputs above.

# --- Row 9 (7. Projektstatusbericht) : brand new reporting row ---
# A9 used to hold the shared formula =1248.94+(2579.69/5); it becomes a plain value.
$ws.Range("A9").Value = 504.25
$ws.Range("B9").Value = 334.13
$ws.Range("C9").Formula = "=A9+B9"
$ws.Range("E9").Value = 490.55
$ws.Range("F9").Value = 334.13
$ws.Range("G9").Formula = "=E9+F9"

# I9/J9/K9 already contain the formulas =A9-E9, =B9-F9, =C9-G9 and now compute
# real numbers instead of the placeholder zeros.

# Row 10 keeps the original formula text (=1248.94+(2579.69/5)) that A9 used to
# share; make sure it stays explicit after A9 became a literal value.
$ws.Range("A10").Formula = "=1248.94+(2579.69/5)"

# Move the active selection as recorded for the saved sheet view.
$ws.Range("D12").Select() | Out-Null

$wb.Application.Calculate()
